# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the cryptos sheet.
# Price cells are forced to text first (NumberFormat "@") so numeric-looking
# values like "298.03" stay stored as strings (matching the source data's
# inline-string cells) instead of being auto-coerced into Excel numbers;
# ClearFormats() afterwards drops the temporary text format again so no
# stray cell style is left behind.
$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.128.38"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.254.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.00%  "
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.78"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.25%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.38"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -7.05%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.603.87"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.268.10"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.065.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("E22").Value = "  -4.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.29"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.76"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.85"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -12.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.55"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.01"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.92"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.33"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("E36").Value = "  -5.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.31"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.67%  "
$ws.Range("E38").Value = "  -5.68%  "
$ws.Range("E39").Value = "  -8.39%  "
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("E41").Value = "  -5.28%  "
$ws.Range("E42").Value = "  -8.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.935.01"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.16%  "
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.21"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.23%  "
$ws.Range("E47").Value = "  -7.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.72%  "
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.482.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.92"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.34%  "
